$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 335
$ws.Range("F3").Value = 1160
$ws.Range("F6").Value = 62
$ws.Range("F9").Value = 1872
$ws.Range("F10").Value = 61
$ws.Range("F11").Value = 458
$ws.Range("F12").Value = 68
$ws.Range("F14").Value = 707
$ws.Range("F15").Value = 475
$ws.Range("F17").Value = 840
$ws.Range("F18").Value = 80740
$ws.Range("F19").Value = 80740
$ws.Range("F21").Value = 689
$ws.Range("F22").Value = 34041
$ws.Range("F23").Value = 34041
$ws.Range("F24").Value = 572
$ws.Range("F27").Value = 68
$ws.Range("F28").Value = 64
$ws.Range("F29").Value = 1032
$ws.Range("F30").Value = 328
$ws.Range("F32").Value = 681
$ws.Range("F33").Value = 3208
$ws.Range("F34").Value = 3208
$ws.Range("F35").Value = 1247
$ws.Range("F36").Value = 5533
$ws.Range("F37").Value = 829
$ws.Range("F38").Value = 482
$ws.Range("F42").Value = 477

$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 1998
$ws.Range("F10").Value = 37
$ws.Range("F14").Value = 12
$ws.Range("F16").Value = 79
$ws.Range("F20").Value = 773
$ws.Range("F39").Value = 37
$ws.Range("F43").Value = 836
$ws.Range("F44").Value = 260

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 737
$ws.Range("F5").Value = 590
$ws.Range("F7").Value = 193

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 335
$ws.Range("F4").Value = 1160
$ws.Range("F10").Value = 62
$ws.Range("F13").Value = 193
$ws.Range("F14").Value = 1872
$ws.Range("F15").Value = 37
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 458
$ws.Range("F18").Value = 68
$ws.Range("F20").Value = 840
$ws.Range("F21").Value = 12
$ws.Range("F23").Value = 80740
$ws.Range("F24").Value = 689
$ws.Range("F25").Value = 34041
$ws.Range("F26").Value = 572
$ws.Range("F31").Value = 64
$ws.Range("F34").Value = 328
$ws.Range("F37").Value = 3208
$ws.Range("F38").Value = 1247
$ws.Range("F39").Value = 5533
$ws.Range("F41").Value = 829
$ws.Range("F47").Value = 477
$ws.Range("F52").Value = 260
